$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to keep text formatting where values look numeric,
# matching the source data (prices are stored as text, not numbers).
# Row 2
$ws.Range("D2").Value = '28.247.31'
$ws.Range("E2").Value = '  -0.94%  '
# Row 3
$ws.Range("D3").Value = '1.806.50'
$ws.Range("E3").Value = '  -1.04%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9969'
$ws.Range("E4").Value = '  -0.74%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.67'
$ws.Range("E5").Value = '  -1.57%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9957'
$ws.Range("E6").Value = '  -0.75%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5156'
$ws.Range("E7").Value = '  -0.22%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3964'
$ws.Range("E8").Value = '  +2.44%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07887'
$ws.Range("E9").Value = '  -6.26%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.106'
$ws.Range("E10").Value = '  -1.29%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.07'
$ws.Range("E11").Value = '  -2.16%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.329'
$ws.Range("E12").Value = '  -1.61%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9962'
$ws.Range("E13").Value = '  -0.73%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.44'
$ws.Range("E14").Value = '  -3.90%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.310'
$ws.Range("E15").Value = '  -2.67%  '
# Row 16
$ws.Range("D16").Value = '1.785.21'
$ws.Range("E16").Value = '  -2.10%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.32'
$ws.Range("E17").Value = '  -2.08%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001081'
$ws.Range("E18").Value = '  -4.51%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06546'
$ws.Range("E19").Value = '  -1.44%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9971'
$ws.Range("E20").Value = '  -0.62%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.26'
$ws.Range("E21").Value = '  -2.92%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.970'
$ws.Range("E22").Value = '  -1.78%  '
# Row 23
$ws.Range("D23").Value = '28.310.92'
$ws.Range("E23").Value = '  -0.88%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -2.72%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.219'
$ws.Range("E25").Value = '  -3.21%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.87'
$ws.Range("E26").Value = '  +0.18%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.56'
$ws.Range("E27").Value = '  -3.43%  '
# Row 28
$ws.Range("D28").Value = '2.003.52'
$ws.Range("E28").Value = '  -1.54%  '
# Row 29
$ws.Range("E29").Value = '  -0.70%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.45'
$ws.Range("E30").Value = '  +1.19%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1081'
$ws.Range("E31").Value = '  -1.35%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  -4.52%  '
# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.578'
$ws.Range("E33").Value = '  -2.92%  '
# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.638'
$ws.Range("E34").Value = '  -1.05%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07137'
$ws.Range("E35").Value = '  -7.95%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.005'
$ws.Range("E36").Value = '  +2.95%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02328'
$ws.Range("E37").Value = '  -2.18%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2144'
$ws.Range("E38").Value = '  -3.79%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.60'
$ws.Range("E39").Value = '  +0.87%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.044'
$ws.Range("E40").Value = '  -4.42%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6180'
$ws.Range("E41").Value = '  -4.01%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9961'
$ws.Range("E42").Value = '  -0.62%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.155'
$ws.Range("E43").Value = '  -3.12%  '
# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  -2.49%  '
# Row 45
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.324'
$ws.Range("E45").Value = '  -5.58%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5975'
$ws.Range("E46").Value = '  -3.52%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.735'
$ws.Range("E47").Value = '  -1.63%  '
# Row 48
$ws.Range("E48").Value = '  -1.66%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.212'
$ws.Range("E49").Value = '  +0.47%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.925'
$ws.Range("E50").Value = '  -3.86%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06867'
$ws.Range("E51").Value = '  -1.84%  '
